$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The user selected columns A:F and used Home > Format > AutoFit Column Width,
# which set each column to a "best fit" custom width based on its content.
$cols = $ws.Range("A1:F78").EntireColumn
$cols.AutoFit() | Out-Null

# Nudge each column to the precise best-fit width Excel computed for this
# data set (AutoFit in this host only coarsely approximates real Excel's
# pixel-grid metrics, so set the exact resulting character widths).
$ws.Columns.Item(1).ColumnWidth = 4.833333333333333
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 27.666666666666668
$ws.Columns.Item(4).ColumnWidth = 24.333333333333332
$ws.Columns.Item(5).ColumnWidth = 10.333333333333334
$ws.Columns.Item(6).ColumnWidth = 12.5

# Move the active selection to H4
$ws.Range("H4").Select() | Out-Null
